$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 48/49 swap: NEARProtocol <-> TrustWalletToken ---
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"

# --- Row 51: HuobiToken -> Celestia ---
$ws.Range("B51").Value = "Celestia"
$ws.Range("C51").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"

# --- Price (column D) updates ---
$ws.Range("D2").Value = "41.989.28"
$ws.Range("D3").Value = "2.235.45"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.629"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "68.49"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.549"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.96"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "34.88"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "2.574.30"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Value = "2.237.98"
$ws.Range("D19").Value = "41.823.62"
$ws.Range("D20").Value = "0.0₃0963"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.126"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0708"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.73"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "21.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.24"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0262"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1000"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.89"
$ws.Range("D51").Style = "Normal"

# --- Volume(1h) (column E) updates ---
$ws.Range("E3").Value = "  -4.98%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("E5").Value = "  -3.81%  "
$ws.Range("E6").Value = "  -5.82%  "
$ws.Range("E7").Value = "  -6.13%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -8.56%  "
$ws.Range("E10").Value = "  -2.78%  "
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("E12").Value = "  +4.97%  "
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("E14").Value = "  -8.04%  "
$ws.Range("E15").Value = "  -4.75%  "
$ws.Range("E16").Value = "  -9.18%  "
$ws.Range("E17").Value = "  -6.23%  "
$ws.Range("E18").Value = "  -4.81%  "
$ws.Range("E20").Value = "  -6.97%  "
$ws.Range("E21").Value = "  -7.66%  "
$ws.Range("E22").Value = "  -7.06%  "
$ws.Range("E23").Value = "  -8.03%  "
$ws.Range("E24").Value = "  +3.31%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  -2.83%  "
$ws.Range("E27").Value = "  -6.42%  "
$ws.Range("E28").Value = "  -6.10%  "
$ws.Range("E29").Value = "  -4.42%  "
$ws.Range("E30").Value = "  -4.76%  "
$ws.Range("E31").Value = "  -9.29%  "
$ws.Range("E32").Value = "  -7.72%  "
$ws.Range("E33").Value = "  -7.98%  "
$ws.Range("E34").Value = "  -2.69%  "
$ws.Range("E35").Value = "  -5.93%  "
$ws.Range("E36").Value = "  -8.01%  "
$ws.Range("E37").Value = "  -5.55%  "
$ws.Range("E38").Value = "  +14.72%  "
$ws.Range("E39").Value = "  -6.02%  "
$ws.Range("E40").Value = "  -7.80%  "
$ws.Range("E41").Value = "  -5.31%  "
$ws.Range("E42").Value = "  -3.15%  "
$ws.Range("E43").Value = "  -3.21%  "
$ws.Range("E44").Value = "  -3.72%  "
$ws.Range("E45").Value = "  -8.79%  "
$ws.Range("E46").Value = "  -6.66%  "
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("E48").Value = "  -6.82%  "
$ws.Range("E49").Value = "  -7.67%  "
$ws.Range("E50").Value = "  +4.17%  "
$ws.Range("E51").Value = "  +3.78%  "
